$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Results")

# --- Step 1: fix up cell formatting (border/wrap styles) for rows whose visual
# style differs from the default, BEFORE overwriting values, by copying formats
# from representative template rows that still carry the old layout. ---

# Template row 19 (RF/version3 originally) carries the wrap+right-border style
# used by the last three rows of each version-4 group; grab it before it changes.
$ws.Range("A19:F19").Copy()
$ws.Range("A25:F27").PasteSpecial(-4122)

# Template row 21 (LightGBM/version3 originally) is the very last data row and has
# no F cell at all; use it for the new last row (29).
$ws.Range("A21:F21").Copy()
$ws.Range("A29:F29").PasteSpecial(-4122)

# Template row 4 (RF/version0) carries the wrap style with the normal left border;
# rows 5 and 6 (the new XGBoost/NeuralNetwork rows for version 0) need this style.
$ws.Range("A4:F4").Copy()
$ws.Range("A5:F6").PasteSpecial(-4122)

# Template row 2 is the plain style (no wrap, normal left border); apply it to every
# other row that must end up in that state.
$ws.Range("A2:F2").Copy()
$ws.Range("A19:F19").PasteSpecial(-4122)
$ws.Range("A21:F21").PasteSpecial(-4122)
$ws.Range("A22:F22").PasteSpecial(-4122)
$ws.Range("A23:F23").PasteSpecial(-4122)
$ws.Range("A24:F24").PasteSpecial(-4122)
$ws.Range("A28:F28").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# --- Step 2: write the final values for every data row (A:E), 7 models x 4 versions ---
$ws.Range("A2").Value2 = 0
$ws.Range("B2").Value2 = "ET"
$ws.Range("C2").Value2 = 0.86832010582010499
$ws.Range("D2").Value2 = 5.99198269844055
$ws.Range("E2").Value2 = 0.75734999999999997

$ws.Range("A3").Value2 = 0
$ws.Range("B3").Value2 = "ET Voting"
$ws.Range("C3").Value2 = 0.86878306878306799
$ws.Range("D3").Value2 = 48.949499368667603
$ws.Range("E3").Value2 = 0.75934000000000001

$ws.Range("A4").Value2 = 0
$ws.Range("B4").Value2 = "RF"
$ws.Range("C4").Value2 = 0.86534391534391497
$ws.Range("D4").Value2 = 6.4758136272430402
$ws.Range("E4").Value2 = 0.75431999999999999

$ws.Range("A5").Value2 = 0
$ws.Range("B5").Value2 = "XGBoost"
$ws.Range("C5").Value2 = 0.84099999999999997
$ws.Range("D5").Value2 = 16.648
$ws.Range("E5").Value2 = 0.79200000000000004

$ws.Range("A6").Value2 = 0
$ws.Range("B6").Value2 = "NeuralNetwork"
$ws.Range("C6").Value2 = 0.95699999999999996
$ws.Range("D6").Value2 = 78.917199999999994
$ws.Range("E6").Value2 = 0.83099999999999996

$ws.Range("A7").Value2 = 0
$ws.Range("B7").Value2 = "Catboost"
$ws.Range("C7").Value2 = 0.85939153439153404
$ws.Range("D7").Value2 = 77.387502193450899
$ws.Range("E7").Value2 = 0.73890999999999996

$ws.Range("A8").Value2 = 0
$ws.Range("B8").Value2 = "LightGBM"
$ws.Range("C8").Value2 = 0.87239999999999995
$ws.Range("D8").Value2 = 56.443899999999999
$ws.Range("E8").Value2 = 0.75253999999999999

$ws.Range("A9").Value2 = 1
$ws.Range("B9").Value2 = "ET"
$ws.Range("C9").Value2 = 0.89093915343915298
$ws.Range("D9").Value2 = 5.7862265110015798
$ws.Range("E9").Value2 = 0.79727000000000003

$ws.Range("A10").Value2 = 1
$ws.Range("B10").Value2 = "ET Voting"
$ws.Range("C10").Value2 = 0.892989417989418
$ws.Range("D10").Value2 = 45.508501768112097
$ws.Range("E10").Value2 = 0.79788999999999999

$ws.Range("A11").Value2 = 1
$ws.Range("B11").Value2 = "RF"
$ws.Range("C11").Value2 = 0.88148148148148098
$ws.Range("D11").Value2 = 6.1565680503845197
$ws.Range("E11").Value2 = 0.78393000000000002

$ws.Range("A12").Value2 = 1
$ws.Range("B12").Value2 = "XGBoost"
$ws.Range("C12").Value2 = 0.85799999999999998
$ws.Range("D12").Value2 = 24.084299999999999
$ws.Range("E12").Value2 = 0.81

$ws.Range("A13").Value2 = 1
$ws.Range("B13").Value2 = "NeuralNetwork"
$ws.Range("C13").Value2 = 0.95299999999999996
$ws.Range("D13").Value2 = 88.8048
$ws.Range("E13").Value2 = 0.84099999999999997

$ws.Range("A14").Value2 = 1
$ws.Range("B14").Value2 = "Catboost"
$ws.Range("C14").Value2 = 0.86216931216931203
$ws.Range("D14").Value2 = 82.544497728347693
$ws.Range("E14").Value2 = 0.74616000000000005

$ws.Range("A15").Value2 = 1
$ws.Range("B15").Value2 = "LightGBM"
$ws.Range("C15").Value2 = 0.88090000000000002
$ws.Range("D15").Value2 = 64.006299999999996
$ws.Range("E15").Value2 = 0.76978000000000002

$ws.Range("A16").Value2 = 2
$ws.Range("B16").Value2 = "ET"
$ws.Range("C16").Value2 = 0.896560846560846
$ws.Range("D16").Value2 = 5.6229577064514098
$ws.Range("E16").Value2 = 0.80054999999999998

$ws.Range("A17").Value2 = 2
$ws.Range("B17").Value2 = "ET Voting"
$ws.Range("C17").Value2 = 0.89755291005290905
$ws.Range("D17").Value2 = 45.323498725891099
$ws.Range("E17").Value2 = 0.80125000000000002

$ws.Range("A18").Value2 = 2
$ws.Range("B18").Value2 = "RF"
$ws.Range("C18").Value2 = 0.88029100529100501
$ws.Range("D18").Value2 = 6.1019423007964999
$ws.Range("E18").Value2 = 0.77829999999999999

$ws.Range("A19").Value2 = 2
$ws.Range("B19").Value2 = "XGBoost"
$ws.Range("C19").Value2 = 0.85899999999999999
$ws.Range("D19").Value2 = 25
$ws.Range("E19").Value2 = 0.80800000000000005

$ws.Range("A20").Value2 = 2
$ws.Range("B20").Value2 = "NeuralNetwork"
$ws.Range("C20").Value2 = 0.96099999999999997
$ws.Range("D20").Value2 = 80.043800000000005
$ws.Range("E20").Value2 = 0.85

$ws.Range("A21").Value2 = 2
$ws.Range("B21").Value2 = "Catboost"
$ws.Range("C21").Value2 = 0.862896825396825
$ws.Range("D21").Value2 = 82.433997392654405
$ws.Range("E21").Value2 = 0.75097999999999998

$ws.Range("A22").Value2 = 2
$ws.Range("B22").Value2 = "LightGBM"
$ws.Range("C22").Value2 = 0.8821
$ws.Range("D22").Value2 = 64.049199999999999
$ws.Range("E22").Value2 = 0.77070000000000005

$ws.Range("A23").Value2 = 3
$ws.Range("B23").Value2 = "ET"
$ws.Range("C23").Value2 = 0.91111111111111098
$ws.Range("D23").Value2 = 5.18176221847534
$ws.Range("E23").Value2 = 0.81398999999999999

$ws.Range("A24").Value2 = 3
$ws.Range("B24").Value2 = "ET Voting"
$ws.Range("C24").Value2 = 0.91309523809523796
$ws.Range("D24").Value2 = 45.2369995117187
$ws.Range("E24").Value2 = 0.81652999999999998

$ws.Range("A25").Value2 = 3
$ws.Range("B25").Value2 = "RF"
$ws.Range("C25").Value2 = 0.89318783068783003
$ws.Range("D25").Value2 = 6.8079984188079798
$ws.Range("E25").Value2 = 0.78535999999999995

$ws.Range("A26").Value2 = 3
$ws.Range("B26").Value2 = "XGBoost"
$ws.Range("C26").Value2 = 0.877
$ws.Range("D26").Value2 = 79
$ws.Range("E26").Value2 = 0.81699999999999995

$ws.Range("A27").Value2 = 3
$ws.Range("B27").Value2 = "NeuralNetwork"
$ws.Range("C27").Value2 = 0.95699999999999996
$ws.Range("D27").Value2 = 79.087699999999998
$ws.Range("E27").Value2 = 0.85299999999999998

$ws.Range("A28").Value2 = 3
$ws.Range("B28").Value2 = "Catboost"
$ws.Range("C28").Value2 = 0.87559523809523798
$ws.Range("D28").Value2 = 96.338500499725299
$ws.Range("E28").Value2 = 0.75775999999999999

$ws.Range("A29").Value2 = 3
$ws.Range("B29").Value2 = "LightGBM"
$ws.Range("C29").Value2 = 0.89670000000000005
$ws.Range("D29").Value2 = 72.608500000000006
$ws.Range("E29").Value2 = 0.78551000000000004

# --- Step 3: grow the Excel Table (ListObject) to cover the new rows ---
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:E29"))

# --- Step 4: restore dimension/sheetView cosmetics ---
$ws.Range("E7").Select()
